$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.477.12"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.902.04"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.19"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.83"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.337"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0702"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0996"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.185.01"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.28"
$ws.Range("E13").Value = "  +7.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.913.36"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.689"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.84"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.530.08"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.88"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0821"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "242.79"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.95"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  +15.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.53"
$ws.Range("E27").Value = "  +7.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.90"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.959"
$ws.Range("E30").Value = "  +24.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0567"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.09"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  +4.65%  "
$ws.Range("E35").Value = "  +6.70%  "
$ws.Range("E36").Value = "  +9.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.02"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.10"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.58"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.73"
$ws.Range("E41").Value = "  +4.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.350.36"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +13.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "49.38"
$ws.Range("E44").Value = "  +42.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.34"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.74"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.65"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.091.49"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0689"
$ws.Range("E51").Value = "  +1.79%  "
